# Case_0_175 loading_percent.xlsx update - "case with 380 kV done"
# Updates computed loading-percent values for rows 2-25, columns C-K and N-O.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 3).Value = 3.685701061400442
$ws.Cells.Item(2, 4).Value = 10.10017929842885
$ws.Cells.Item(2, 5).Value = 14.03881895929419
$ws.Cells.Item(2, 6).Value = 28.29355485500399
$ws.Cells.Item(2, 7).Value = 26.86987343044801
$ws.Cells.Item(2, 8).Value = 13.78051638035478
$ws.Cells.Item(2, 9).Value = 19.40484108690837
$ws.Cells.Item(2, 10).Value = 9.740175377643377
$ws.Cells.Item(2, 11).Value = 13.36109662984896
$ws.Cells.Item(2, 14).Value = 16.65595541911686
$ws.Cells.Item(2, 15).Value = 20.73030004301128

# Row 3
$ws.Cells.Item(3, 3).Value = 3.555038748699812
$ws.Cells.Item(3, 4).Value = 10.04426373865938
$ws.Cells.Item(3, 5).Value = 13.97484461318791
$ws.Cells.Item(3, 6).Value = 28.3126419055888
$ws.Cells.Item(3, 7).Value = 26.87990778699228
$ws.Cells.Item(3, 8).Value = 13.82506560225322
$ws.Cells.Item(3, 9).Value = 19.44354844778263
$ws.Cells.Item(3, 10).Value = 9.745206472895356
$ws.Cells.Item(3, 11).Value = 12.86315344576962
$ws.Cells.Item(3, 14).Value = 16.6673787505694
$ws.Cells.Item(3, 15).Value = 20.79188227735693

# Row 4
$ws.Cells.Item(4, 3).Value = 3.471517390418423
$ws.Cells.Item(4, 4).Value = 10.01147868663936
$ws.Cells.Item(4, 5).Value = 13.9383277492021
$ws.Cells.Item(4, 6).Value = 28.33204102417042
$ws.Cells.Item(4, 7).Value = 26.89609832379393
$ws.Cells.Item(4, 8).Value = 13.8548733790866
$ws.Cells.Item(4, 9).Value = 19.47226260555254
$ws.Cells.Item(4, 10).Value = 9.749864171428282
$ws.Cells.Item(4, 11).Value = 12.5478718223984
$ws.Cells.Item(4, 14).Value = 16.67607931596168
$ws.Cells.Item(4, 15).Value = 20.83476068856033

# Row 5
$ws.Cells.Item(5, 3).Value = 3.436682564982792
$ws.Cells.Item(5, 4).Value = 9.998518230782578
$ws.Cells.Item(5, 5).Value = 13.9241536189655
$ws.Cells.Item(5, 6).Value = 28.34187496469368
$ws.Cells.Item(5, 7).Value = 26.9052112047065
$ws.Cells.Item(5, 8).Value = 13.86763688998082
$ws.Cells.Item(5, 9).Value = 19.48520496011356
$ws.Cells.Item(5, 10).Value = 9.752157074199083
$ws.Cells.Item(5, 11).Value = 12.41717038944968
$ws.Cells.Item(5, 14).Value = 16.68005014164614
$ws.Cells.Item(5, 15).Value = 20.8535043427433

# Row 6
$ws.Cells.Item(6, 3).Value = 3.430850918680331
$ws.Cells.Item(6, 4).Value = 9.996390598755243
$ws.Cells.Item(6, 5).Value = 13.9218430316538
$ws.Cells.Item(6, 6).Value = 28.34362427535189
$ws.Cells.Item(6, 7).Value = 26.90687605393071
$ws.Cells.Item(6, 8).Value = 13.86979348536932
$ws.Cells.Item(6, 9).Value = 19.48742889145481
$ws.Cells.Item(6, 10).Value = 9.752561667035419
$ws.Cells.Item(6, 11).Value = 12.39533891663624
$ws.Cells.Item(6, 14).Value = 16.68073521316349
$ws.Cells.Item(6, 15).Value = 20.85669332621804

# Row 7
$ws.Cells.Item(7, 3).Value = 3.471050790662968
$ws.Cells.Item(7, 4).Value = 10.01130226531676
$ws.Cells.Item(7, 5).Value = 13.93813371519946
$ws.Cells.Item(7, 6).Value = 28.33216584319461
$ws.Cells.Item(7, 7).Value = 26.8962110504606
$ws.Cells.Item(7, 8).Value = 13.85504301644284
$ws.Cells.Item(7, 9).Value = 19.47243212980484
$ws.Cells.Item(7, 10).Value = 9.749893495105296
$ws.Cells.Item(7, 11).Value = 12.54611788306706
$ws.Cells.Item(7, 14).Value = 16.67613114447325
$ws.Cells.Item(7, 15).Value = 20.83500833348171

# Row 8
$ws.Cells.Item(8, 3).Value = 3.641348845745696
$ws.Cells.Item(8, 4).Value = 10.08058500705021
$ws.Cells.Item(8, 5).Value = 14.01619439541276
$ws.Cells.Item(8, 6).Value = 28.2985410137017
$ws.Cells.Item(8, 7).Value = 26.87124821144263
$ws.Cells.Item(8, 8).Value = 13.79536720746356
$ws.Cells.Item(8, 9).Value = 19.41715878461884
$ws.Cells.Item(8, 10).Value = 9.741584823382411
$ws.Cells.Item(8, 11).Value = 13.19148718130634
$ws.Cells.Item(8, 14).Value = 16.65954484001583
$ws.Cells.Item(8, 15).Value = 20.75047984504634

# Row 9
$ws.Cells.Item(9, 3).Value = 3.948052320803114
$ws.Cells.Item(9, 4).Value = 10.22823446162508
$ws.Cells.Item(9, 5).Value = 14.19063741710212
$ws.Cells.Item(9, 6).Value = 28.29360079965121
$ws.Cells.Item(9, 7).Value = 26.90207821928456
$ws.Cells.Item(9, 8).Value = 13.69784266238242
$ws.Cells.Item(9, 9).Value = 19.34814516303977
$ws.Cells.Item(9, 10).Value = 9.737715377217985
$ws.Cells.Item(9, 11).Value = 14.37447854956842
$ws.Cells.Item(9, 14).Value = 16.64035011760251
$ws.Cells.Item(9, 15).Value = 20.62507882433599

# Row 10
$ws.Cells.Item(10, 3).Value = 4.155500160014094
$ws.Cells.Item(10, 4).Value = 10.3432081230394
$ws.Cells.Item(10, 5).Value = 14.33101686431509
$ws.Cells.Item(10, 6).Value = 28.32717684087627
$ws.Cells.Item(10, 7).Value = 26.9734909547842
$ws.Cells.Item(10, 8).Value = 13.63811258795734
$ws.Cells.Item(10, 9).Value = 19.32158569704688
$ws.Cells.Item(10, 10).Value = 9.742410961540198
$ws.Cells.Item(10, 11).Value = 15.18534526996618
$ws.Cells.Item(10, 14).Value = 16.63430537714756
$ws.Cells.Item(10, 15).Value = 20.55775329415963

# Row 11
$ws.Cells.Item(11, 3).Value = 4.245774987704022
$ws.Cells.Item(11, 4).Value = 10.39676764744955
$ws.Cells.Item(11, 5).Value = 14.39734322658498
$ws.Cells.Item(11, 6).Value = 28.35051004596015
$ws.Cells.Item(11, 7).Value = 27.01654256564369
$ws.Cells.Item(11, 8).Value = 13.61353538055983
$ws.Cells.Item(11, 9).Value = 19.314766593829
$ws.Cells.Item(11, 10).Value = 9.746173923035711
$ws.Cells.Item(11, 11).Value = 15.54022419891396
$ws.Cells.Item(11, 14).Value = 16.63328959964158
$ws.Cells.Item(11, 15).Value = 20.53255009706135

# Row 12
$ws.Cells.Item(12, 3).Value = 4.279356168401943
$ws.Cells.Item(12, 4).Value = 10.41721644995536
$ws.Cells.Item(12, 5).Value = 14.42279712143154
$ws.Cells.Item(12, 6).Value = 28.36050082341885
$ws.Cells.Item(12, 7).Value = 27.03435833482185
$ws.Cells.Item(12, 8).Value = 13.60460226225663
$ws.Cells.Item(12, 9).Value = 19.31294204425344
$ws.Cells.Item(12, 10).Value = 9.747831705605927
$ws.Cells.Item(12, 11).Value = 15.67250269638419
$ws.Cells.Item(12, 14).Value = 16.63315289036331
$ws.Cells.Item(12, 15).Value = 20.52378890636884

# Row 13
$ws.Cells.Item(13, 3).Value = 4.272150940100531
$ws.Cells.Item(13, 4).Value = 10.41280523139138
$ws.Cells.Item(13, 5).Value = 14.41730043776852
$ws.Cells.Item(13, 6).Value = 28.35829783751238
$ws.Cells.Item(13, 7).Value = 27.03045422283605
$ws.Cells.Item(13, 8).Value = 13.60650953650477
$ws.Cells.Item(13, 9).Value = 19.31330128601655
$ws.Cells.Item(13, 10).Value = 9.747464336501894
$ws.Cells.Item(13, 11).Value = 15.64410921615953
$ws.Cells.Item(13, 14).Value = 16.63317132866761
$ws.Cells.Item(13, 15).Value = 20.52564093244952

# Row 14
$ws.Cells.Item(14, 3).Value = 4.248549914971567
$ws.Cells.Item(14, 4).Value = 10.39844671654559
$ws.Cells.Item(14, 5).Value = 14.39943066967554
$ws.Cells.Item(14, 6).Value = 28.35130887221404
$ws.Cells.Item(14, 7).Value = 27.017978000525
$ws.Cells.Item(14, 8).Value = 13.6127929534915
$ws.Cells.Item(14, 9).Value = 19.31460129833748
$ws.Cells.Item(14, 10).Value = 9.746305650415644
$ws.Cells.Item(14, 11).Value = 15.55114954521689
$ws.Cells.Item(14, 14).Value = 16.63327339154289
$ws.Cells.Item(14, 15).Value = 20.53181360680632

# Row 15
$ws.Cells.Item(15, 3).Value = 4.234014540436621
$ws.Cells.Item(15, 4).Value = 10.3896730411941
$ws.Cells.Item(15, 5).Value = 14.38852835765539
$ws.Cells.Item(15, 6).Value = 28.34717820654319
$ws.Cells.Item(15, 7).Value = 27.01053277036491
$ws.Cells.Item(15, 8).Value = 13.61669041914712
$ws.Cells.Item(15, 9).Value = 19.31549628584125
$ws.Cells.Item(15, 10).Value = 9.745626207674754
$ws.Cells.Item(15, 11).Value = 15.49393214227224
$ws.Cells.Item(15, 14).Value = 16.63336815507819
$ws.Cells.Item(15, 15).Value = 20.53569655901496

# Row 16
$ws.Cells.Item(16, 3).Value = 4.149516753667957
$ws.Cells.Item(16, 4).Value = 10.33973200748791
$ws.Cells.Item(16, 5).Value = 14.32673046555628
$ws.Cells.Item(16, 6).Value = 28.32581382830939
$ws.Cells.Item(16, 7).Value = 26.97088955966754
$ws.Cells.Item(16, 8).Value = 13.63977102883925
$ws.Cells.Item(16, 9).Value = 19.32213733297296
$ws.Cells.Item(16, 10).Value = 9.742197688410513
$ws.Cells.Item(16, 11).Value = 15.16186320245052
$ws.Cells.Item(16, 14).Value = 16.63440653428115
$ws.Cells.Item(16, 15).Value = 20.55950976499832

# Row 17
$ws.Cells.Item(17, 3).Value = 4.09661996952001
$ws.Cells.Item(17, 4).Value = 10.30940709513593
$ws.Cells.Item(17, 5).Value = 14.28943927086456
$ws.Cells.Item(17, 6).Value = 28.31476926239628
$ws.Cells.Item(17, 7).Value = 26.949272179552
$ws.Cells.Item(17, 8).Value = 13.65459518886693
$ws.Cells.Item(17, 9).Value = 19.32756015810816
$ws.Cells.Item(17, 10).Value = 9.740510342881617
$ws.Cells.Item(17, 11).Value = 14.95449517573248
$ws.Cells.Item(17, 14).Value = 16.63548668284519
$ws.Cells.Item(17, 15).Value = 20.57550958672577

# Row 18
$ws.Cells.Item(18, 3).Value = 4.065810533231766
$ws.Cells.Item(18, 4).Value = 10.29208421000926
$ws.Cells.Item(18, 5).Value = 14.26822364273664
$ws.Cells.Item(18, 6).Value = 28.30917561294618
$ws.Cells.Item(18, 7).Value = 26.93783345315928
$ws.Cells.Item(18, 8).Value = 13.66336580244399
$ws.Cells.Item(18, 9).Value = 19.33117457369793
$ws.Cells.Item(18, 10).Value = 9.73969303492939
$ws.Cells.Item(18, 11).Value = 14.83391071649655
$ws.Cells.Item(18, 14).Value = 16.63627125552191
$ws.Cells.Item(18, 15).Value = 20.58522259837725

# Row 19
$ws.Cells.Item(19, 3).Value = 4.055313429401749
$ws.Cells.Item(19, 4).Value = 10.28623986929205
$ws.Cells.Item(19, 5).Value = 14.26108097661877
$ws.Cells.Item(19, 6).Value = 28.30741213884377
$ws.Cells.Item(19, 7).Value = 26.93413154209264
$ws.Cells.Item(19, 8).Value = 13.66637729324818
$ws.Cells.Item(19, 9).Value = 19.33248339412484
$ws.Cells.Item(19, 10).Value = 9.739442654062863
$ws.Cells.Item(19, 11).Value = 14.79286068514261
$ws.Cells.Item(19, 14).Value = 16.63656498724786
$ws.Cells.Item(19, 15).Value = 20.58859881758945

# Row 20
$ws.Cells.Item(20, 3).Value = 4.102290844480275
$ws.Cells.Item(20, 4).Value = 10.31262299627536
$ws.Cells.Item(20, 5).Value = 14.29338496137807
$ws.Cells.Item(20, 6).Value = 28.31586645713888
$ws.Cells.Item(20, 7).Value = 26.95147044322013
$ws.Cells.Item(20, 8).Value = 13.65299185823902
$ws.Cells.Item(20, 9).Value = 19.32693161416962
$ws.Cells.Item(20, 10).Value = 9.740674114851137
$ws.Cells.Item(20, 11).Value = 14.97670635695581
$ws.Cells.Item(20, 14).Value = 16.63535480828561
$ws.Cells.Item(20, 15).Value = 20.57375353903908

# Row 21
$ws.Cells.Item(21, 3).Value = 4.255498609443137
$ws.Cells.Item(21, 4).Value = 10.40265974039481
$ws.Cells.Item(21, 5).Value = 14.40467043303888
$ws.Cells.Item(21, 6).Value = 28.35333039180911
$ws.Cells.Item(21, 7).Value = 27.02160156779844
$ws.Cells.Item(21, 8).Value = 13.61093721320342
$ws.Cells.Item(21, 9).Value = 19.31419888490821
$ws.Cells.Item(21, 10).Value = 9.746639675227998
$ws.Cells.Item(21, 11).Value = 15.5785119155546
$ws.Cells.Item(21, 14).Value = 16.63323669551453
$ws.Cells.Item(21, 15).Value = 20.52997927912116

# Row 22
$ws.Cells.Item(22, 3).Value = 4.352102929183775
$ws.Cells.Item(22, 4).Value = 10.46247061645757
$ws.Cells.Item(22, 5).Value = 14.47936099315682
$ws.Cells.Item(22, 6).Value = 28.38454475862011
$ws.Cells.Item(22, 7).Value = 27.07625090623855
$ws.Cells.Item(22, 8).Value = 13.58563079391699
$ws.Cells.Item(22, 9).Value = 19.31029380373531
$ws.Cells.Item(22, 10).Value = 9.751895015214314
$ws.Cells.Item(22, 11).Value = 15.95951335832378
$ws.Cells.Item(22, 14).Value = 16.63329690581477
$ws.Cells.Item(22, 15).Value = 20.50593361710156

# Row 23
$ws.Cells.Item(23, 3).Value = 4.30087048086978
$ws.Cells.Item(23, 4).Value = 10.4304646129511
$ws.Cells.Item(23, 5).Value = 14.43932381332985
$ws.Cells.Item(23, 6).Value = 28.36727092909627
$ws.Cells.Item(23, 7).Value = 27.04627970482785
$ws.Cells.Item(23, 8).Value = 13.59893773454851
$ws.Cells.Item(23, 9).Value = 19.31197374509717
$ws.Cells.Item(23, 10).Value = 9.748966425698468
$ws.Cells.Item(23, 11).Value = 15.75732049767745
$ws.Cells.Item(23, 14).Value = 16.63313307326752
$ws.Cells.Item(23, 15).Value = 20.5183488388548

# Row 24
$ws.Cells.Item(24, 3).Value = 4.099728281583252
$ws.Cells.Item(24, 4).Value = 10.31116873960661
$ws.Cells.Item(24, 5).Value = 14.2916004169913
$ws.Cells.Item(24, 6).Value = 28.31536806016249
$ws.Cells.Item(24, 7).Value = 26.95047352587537
$ws.Cells.Item(24, 8).Value = 13.65371595153134
$ws.Cells.Item(24, 9).Value = 19.32721423157476
$ws.Cells.Item(24, 10).Value = 9.740599597600925
$ws.Cells.Item(24, 11).Value = 14.96666892879299
$ws.Cells.Item(24, 14).Value = 16.63541391912127
$ws.Cells.Item(24, 15).Value = 20.57454584576454

# Row 25
$ws.Cells.Item(25, 3).Value = 3.868140856973458
$ws.Cells.Item(25, 4).Value = 10.18709901163442
$ws.Cells.Item(25, 5).Value = 14.14123888981922
$ws.Cells.Item(25, 6).Value = 28.28840028789804
$ws.Cells.Item(25, 7).Value = 26.88517425668568
$ws.Cells.Item(25, 8).Value = 13.72213452363061
$ws.Cells.Item(25, 9).Value = 19.36258340786074
$ws.Cells.Item(25, 10).Value = 9.737435254573917
$ws.Cells.Item(25, 11).Value = 14.06417315194522
$ws.Cells.Item(25, 14).Value = 16.64412260557595
$ws.Cells.Item(25, 15).Value = 20.65466051443924
